$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 221.0990153490724
$ws.Range("R2").Value = 1989.891138141652
$ws.Range("S2").Value = 0.001818510806315634
$ws.Range("T2").Value = 0.001818510806315634
$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 20888.29669366387
$ws.Range("R3").Value = 187994.6702429748
$ws.Range("S3").Value = 0.1718035388035672
$ws.Range("T3").Value = 0.1718035388035672
$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 3346.912582192649
$ws.Range("R4").Value = 30122.21323973384
$ws.Range("S4").Value = 0.02752792312938099
$ws.Range("T4").Value = 0.02752792312938099
$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 9433.69261799597
$ws.Range("R5").Value = 84903.23356196373
$ws.Range("S5").Value = 0.07759090171523764
$ws.Range("T5").Value = 0.07759090171523762
$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 354.694104522546
$ws.Range("R6").Value = 3192.246940702913
$ws.Range("S6").Value = 0.002917313136796848
$ws.Range("T6").Value = 0.002917313136796848
$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.275612726060887
$ws.Range("T7").Value = 0.275612726060887
$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 5369.224097998878
$ws.Range("R8").Value = 48323.0168819899
$ws.Range("S8").Value = 0.04416117379955684
$ws.Range("T8").Value = 0.04416117379955684
$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 15133.83110367194
$ws.Range("R9").Value = 136204.4799330475
$ws.Range("S9").Value = 0.1244738035559892
$ws.Range("T9").Value = 0.1244738035559892
$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 116.1916142370589
$ws.Range("R10").Value = 1045.72452813353
$ws.Range("S10").Value = 0.0009556609999359526
$ws.Range("T10").Value = 0.0009556609999359524
$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 10977.18552779438
$ws.Range("R11").Value = 98794.66975014939
$ws.Range("S11").Value = 0.0902859312770314
$ws.Range("T11").Value = 0.09028593127703138
$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 1758.864348723199
$ws.Range("R12").Value = 15829.77913850879
$ws.Range("S12").Value = 0.01446643179277219
$ws.Range("T12").Value = 0.01446643179277219
$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 4957.579624543436
$ws.Range("R13").Value = 44618.21662089092
$ws.Range("S13").Value = 0.04077545124372832
$ws.Range("T13").Value = 0.04077545124372831
$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 101.2215575116861
$ws.Range("R14").Value = 910.994017605175
$ws.Range("S14").Value = 0.0008325342194603892
$ws.Range("T14").Value = 0.0008325342194603892
$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 9562.891638213396
$ws.Range("R15").Value = 86066.02474392056
$ws.Range("S15").Value = 0.07865354694710282
$ws.Range("T15").Value = 0.07865354694710282
$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 1532.253338578336
$ws.Range("R16").Value = 13790.28004720503
$ws.Range("S16").Value = 0.01260258554213232
$ws.Range("T16").Value = 0.01260258554213232
$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 4318.848088818744
$ws.Range("R17").Value = 38869.6327993687
$ws.Range("S17").Value = 0.0355219669701051
$ws.Range("T17").Value = 0.0355219669701051
